# "Document for presentation and update excel"
#
# Adds a new tracked time-entry row ("review before presentation") right
# before the totals row, and refreshes the totals row so it sums through
# the new row and reports the new total duration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing totals row (old row 17) down to make room for the new
# time entry; the new row 17 inherits row 16's formatting (date / time /
# text / minutes styles), same as it would if the author had inserted the
# row above the totals in Excel.
$null = $ws.Rows.Item(17).Insert()

# --- refresh the totals row (now row 18) -------------------------------
$ws.Range("E18").Value = "total (11 hours and 25 minutes)"
$ws.Range("F18").Formula = "=SUM(F6:F17)"

# --- fill in the new time entry (row 17) -------------------------------
$ws.Range("B17").Value = 45165
$ws.Range("C17").Value = 0.95833333333333337
$ws.Range("D17").Value = 0.97916666666666663
$ws.Range("E17").Value = "review before presentation"
$ws.Range("F17").Formula = "=(D17-C17)*60*24"

# Keep the active selection where the author left it after the edit.
$null = $ws.Range("E18").Select()
